# Add a new "time_taken" column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same text style as the other header cells (B1:E1).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells F2:F23: per-row timestamp strings.
$ws.Range("F2").Value = "2021-10-05 10:50:44.317433"
$ws.Range("F3").Value = "2021-10-05 10:50:44.317445"
$ws.Range("F4").Value = "2021-10-05 10:50:44.317448"
$ws.Range("F5").Value = "2021-10-05 10:50:44.317451"
$ws.Range("F6").Value = "2021-10-05 10:50:44.317454"
$ws.Range("F7").Value = "2021-10-05 10:50:44.317457"
$ws.Range("F8").Value = "2021-10-05 10:50:44.317460"
$ws.Range("F9").Value = "2021-10-05 10:50:44.317463"
$ws.Range("F10").Value = "2021-10-05 10:50:44.317466"
$ws.Range("F11").Value = "2021-10-05 10:50:44.317469"
$ws.Range("F12").Value = "2021-10-05 10:50:44.317471"
$ws.Range("F13").Value = "2021-10-05 10:50:44.317474"
$ws.Range("F14").Value = "2021-10-05 10:50:44.317476"
$ws.Range("F15").Value = "2021-10-05 10:50:44.317479"
$ws.Range("F16").Value = "2021-10-05 10:50:44.317482"
$ws.Range("F17").Value = "2021-10-05 10:50:44.317486"
$ws.Range("F18").Value = "2021-10-05 10:50:44.317491"
$ws.Range("F19").Value = "2021-10-05 10:50:44.317495"
$ws.Range("F20").Value = "2021-10-05 10:50:44.317500"
$ws.Range("F21").Value = "2021-10-05 10:50:44.317504"
$ws.Range("F22").Value = "2021-10-05 10:50:44.317508"
$ws.Range("F23").Value = "2021-10-05 10:50:44.317513"
